$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need special handling so
# Excel keeps storing them as text (matching the original inline-string cells)
# instead of silently converting them to a numeric value.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

$ws.Range('D2').Value = '66.987.27'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.594.41'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '588.84'
$ws.Range('E5').Value = '  -1.99%  '
Set-TextValue 'D6' '149.23'
$ws.Range('E6').Value = '  -3.43%  '
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').Value = '2.592.29'
$ws.Range('E9').Value = '  -0.94%  '
Set-TextValue 'D10' '0.123'
$ws.Range('E10').Value = '  -3.25%  '
$ws.Range('E12').Value = '  -1.91%  '
Set-TextValue 'D13' '0.343'
$ws.Range('E13').Value = '  -3.12%  '
Set-TextValue 'D14' '27.09'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '3.062.43'
$ws.Range('E15').Value = '  -1.10%  '
Set-TextValue 'D16' '0.0000178'
$ws.Range('E16').Value = '  -4.91%  '
$ws.Range('D17').Value = '66.944.72'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '2.593.30'
$ws.Range('E18').Value = '  -0.81%  '
Set-TextValue 'D19' '362.42'
$ws.Range('E19').Value = '  -1.00%  '
Set-TextValue 'D20' '10.96'
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('E21').Value = '  -4.54%  '
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('E24').Value = '  -0.19%  '
Set-TextValue 'D25' '72.20'
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('E28').Value = '  -1.24%  '
Set-TextValue 'D30' '574.78'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('D31').Value = '0.0₃0976'
$ws.Range('E31').Value = '  -6.18%  '
$ws.Range('E32').Value = '  -4.99%  '
Set-TextValue 'D33' '7.58'
$ws.Range('E33').Value = '  -4.05%  '
Set-TextValue 'D34' '1.80'
$ws.Range('E34').Value = '  -3.58%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -5.98%  '
$ws.Range('E37').Value = '  -2.97%  '
Set-TextValue 'D38' '156.25'
$ws.Range('E38').Value = '  -1.23%  '
Set-TextValue 'D39' '18.87'
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('E40').Value = '  -1.46%  '
$ws.Range('E41').Value = '  -1.20%  '
Set-TextValue 'D42' '5.17'
$ws.Range('E42').Value = '  -3.33%  '
Set-TextValue 'D44' '2.47'
$ws.Range('E44').Value = '  -4.76%  '
$ws.Range('E45').Value = '  -0.06%  '
Set-TextValue 'D46' '151.65'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('D47').Value = '0.0₆0281'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('E49').Value = '  -2.81%  '
Set-TextValue 'D50' '0.0775'
$ws.Range('E50').Value = '  -1.85%  '
Set-TextValue 'D51' '21.19'
$ws.Range('E51').Value = '  +1.15%  '
